$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnes")

$names = @("jack", "queen", "king", "carnis")
$startId = 5
$startRow = 6

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $startId + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).NumberFormat = "General"
}
